# Update LDLC prices history
# Insert a new timestamp snapshot column ("EY") right before the
# "nom" / "url_produit" columns, shifting them one column to the right
# (old EY -> EZ, old EZ -> FA), and fill the new column with the latest
# known price (copied from the previous snapshot column, now EX) for
# every product row that still has price data, leaving it blank for
# rows whose tracking already stopped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column EY is column number 155 (A=1). Inserting here pushes the
# former "nom" (EY) and "url_produit" (EZ) columns out to EZ and FA,
# exactly matching the diff's dimension change A1:EZ206 -> A1:FA206.
$newColIndex = 155
$lastSnapshotColIndex = 154

$ws.Columns.Item($newColIndex).Insert()

# Header cell for the newly inserted snapshot column.
$ws.Cells.Item(1, $newColIndex).Value2 = "2026-02-04 00:58:04"

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $prevVal = $ws.Cells.Item($r, $lastSnapshotColIndex).Value2
    if ($prevVal -eq "") {
        $ws.Cells.Item($r, $newColIndex).Value2 = ""
    } else {
        $ws.Cells.Item($r, $newColIndex).Value2 = $prevVal
    }
}
